$wb = $excel.ActiveWorkbook

$wsGlm2 = $wb.Worksheets.Item("glm2")
$wsBeta = $wb.Worksheets.Item("beta")

# ---------------------------------------------------------------------
# glm2 sheet ("B" column = new fitted coefficients for the logit-link
# GLM model; column C onward keeps the previous / reference values).
# Setting Font.Color on these cells (explicit black) mirrors the style
# bump from the default / exponential style to the "no-format" (s=8) or
# "exponential" (s=9) black-font styles already used elsewhere in the
# workbook.
# ---------------------------------------------------------------------
$wsGlm2.Cells.Item(2, 2).Value = -5.301350553
$wsGlm2.Cells.Item(2, 2).Font.Color = 0

$wsGlm2.Cells.Item(3, 2).Value = 7.265002728
$wsGlm2.Cells.Item(3, 2).Font.Color = 0

$wsGlm2.Cells.Item(4, 2).Value = 0.440989745
$wsGlm2.Cells.Item(4, 2).Font.Color = 0

$wsGlm2.Cells.Item(5, 2).Value = 0.000232767
$wsGlm2.Cells.Item(5, 2).Font.Color = 0

$wsGlm2.Cells.Item(6, 2).Value = -0.00067813
$wsGlm2.Cells.Item(6, 2).Font.Color = 0

$wsGlm2.Cells.Item(7, 2).Value = -0.003050189
$wsGlm2.Cells.Item(7, 2).Font.Color = 0

$wsGlm2.Cells.Item(8, 2).Value = 0.001177158
$wsGlm2.Cells.Item(8, 2).Font.Color = 0

$wsGlm2.Cells.Item(9, 2).Value = -0.013657028
$wsGlm2.Cells.Item(9, 2).Font.Color = 0

$wsGlm2.Cells.Item(10, 2).Value = 0.0000324
$wsGlm2.Cells.Item(10, 2).Font.Color = 0

$wsGlm2.Cells.Item(11, 2).Value = -0.025872469
$wsGlm2.Cells.Item(11, 2).Font.Color = 0

$wsGlm2.Cells.Item(12, 2).Value = -0.627988058
$wsGlm2.Cells.Item(12, 2).Font.Color = 0

$wsGlm2.Cells.Item(13, 2).Value = -0.181555015
$wsGlm2.Cells.Item(13, 2).Font.Color = 0

$wsGlm2.Cells.Item(14, 2).Value = -0.036474089
$wsGlm2.Cells.Item(14, 2).Font.Color = 0

$wsGlm2.Cells.Item(15, 2).Value = -0.021084725
$wsGlm2.Cells.Item(15, 2).Font.Color = 0

$wsGlm2.Cells.Item(16, 2).Value = -0.002131749
$wsGlm2.Cells.Item(16, 2).Font.Color = 0

# ---------------------------------------------------------------------
# beta sheet: column B holds the new "beta" coefficients, columns C:G
# are refreshed to the (now updated) glm2 "B" column value for that
# variable (they mirror / broadcast the glm2 result across C:G).
# ---------------------------------------------------------------------
function Set-BetaRow($row, $bValue, $cgValue) {
    $wsBeta.Cells.Item($row, 2).Value = $bValue
    $wsBeta.Cells.Item($row, 3).Value = $cgValue
    $wsBeta.Cells.Item($row, 4).Value = $cgValue
    $wsBeta.Cells.Item($row, 5).Value = $cgValue
    $wsBeta.Cells.Item($row, 6).Value = $cgValue
    $wsBeta.Cells.Item($row, 7).Value = $cgValue
}

Set-BetaRow 2  -5.033172522  -5.301350553
Set-BetaRow 3  9.979031699   7.265002728
Set-BetaRow 4  0.26250815    0.440989745
Set-BetaRow 5  0.000365907   0.000232767
Set-BetaRow 6  -0.000651392  -0.00067813
Set-BetaRow 7  -0.000840751  -0.003050189
Set-BetaRow 8  -0.000320507  0.001177158
Set-BetaRow 9  -0.002275854  -0.013657028
Set-BetaRow 10 0.00000636    0.0000324
Set-BetaRow 11 -0.014953316  -0.025872469
Set-BetaRow 12 -0.250141893  -0.627988058
Set-BetaRow 13 -0.104370497  -0.181555015
Set-BetaRow 14 -0.047374274  -0.036474089
Set-BetaRow 15 -0.025258573  -0.021084725
Set-BetaRow 16 -0.024602938  -0.002131749

# rows 5 and 6 lose their exponential number format (s=9 -> s=8) while
# row 10 keeps it (s=9) -- set explicit black font on every row so the
# font-based style (8/9) is applied consistently, then restore the
# exponential number format on row 10 which must stay s=9.
$wsBeta.Cells.Item(5, 2).Font.Color = 0
$wsBeta.Cells.Item(5, 3).Font.Color = 0
$wsBeta.Cells.Item(5, 4).Font.Color = 0
$wsBeta.Cells.Item(5, 5).Font.Color = 0
$wsBeta.Cells.Item(5, 6).Font.Color = 0
$wsBeta.Cells.Item(5, 7).Font.Color = 0

$wsBeta.Cells.Item(6, 2).Font.Color = 0
$wsBeta.Cells.Item(6, 3).Font.Color = 0
$wsBeta.Cells.Item(6, 4).Font.Color = 0
$wsBeta.Cells.Item(6, 5).Font.Color = 0
$wsBeta.Cells.Item(6, 6).Font.Color = 0
$wsBeta.Cells.Item(6, 7).Font.Color = 0

# ---------------------------------------------------------------------
# beta sheet: the old trailing blank row (17) is no longer needed now
# that data stops at row 16.
# ---------------------------------------------------------------------
$wsBeta.Rows.Item(17).Delete()

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping: glm2 becomes the active tab
# with B8 selected; beta is left scrolled to, and selecting, the old
# (now-empty) row-17-and-below region.
# ---------------------------------------------------------------------
$wsBeta.Activate()
$wsBeta.Range("A17:XFD46").Select()

$wsGlm2.Activate()
$wsGlm2.Range("B8").Select()
